# Auto-generated Excel COM-interop edit script
# Applies the weekly CompStat data refresh described by the commit "New crime data collected".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Volume/Number and "Report Covering the Week" date range ---
# "Volume 31   Number  51" -> "...  52" (characters 21-22 are the "51")
$ws.Range("A8").Characters(21, 2).Text = "52"

# "Report Covering the Week  12/16/2024  Through  12/22/2024"
#  -> "...  12/23/2024  Through  12/29/2024"
$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "12/23/2024"
$c9.Characters(48, 10).Text = "12/29/2024"

# --- Crime-statistics grid (rows 14-31): refreshed weekly/28-day/YTD/2-year figures ---

# A few cells flip between the literal placeholder text ("0" / "***.*") and real
# numbers this week, so their number format is pinned explicitly before the value
# is written (matches how Excel itself decides text vs. number storage).

# Row 14
$ws.Range("L14").Value = -50
$ws.Range("N14").Value = -63.636363636363
# Row 15
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 30
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 30.434782608695
$ws.Range("M15").Value = 76.470588235294
$ws.Range("N15").Value = -25
# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 401
$ws.Range("J16").Value = 445
$ws.Range("K16").Value = -9.887640449438
$ws.Range("L16").Value = -37.732919254658
$ws.Range("M16").Value = 138.690476190476
$ws.Range("N16").Value = -83.94715772618
# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -9.090909090909
$ws.Range("I17").Value = 530
$ws.Range("J17").Value = 479
$ws.Range("K17").Value = 10.647181628392
$ws.Range("L17").Value = 7.505070993914
$ws.Range("M17").Value = 191.208791208791
$ws.Range("N17").Value = -20.180722891566
# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -13.043478260869
$ws.Range("I18").Value = 343
$ws.Range("J18").Value = 406
$ws.Range("K18").Value = -15.51724137931
$ws.Range("L18").Value = -47.872340425531
$ws.Range("M18").Value = -3.38028169014
$ws.Range("N18").Value = -87.143928035982
# Row 19
$ws.Range("C19").Value = 45
$ws.Range("D19").Value = 57
$ws.Range("E19").Value = -21.052631578947
$ws.Range("F19").Value = 138
$ws.Range("G19").Value = 223
$ws.Range("H19").Value = -38.116591928251
$ws.Range("I19").Value = 2036
$ws.Range("J19").Value = 2371
$ws.Range("K19").Value = -14.129059468578
$ws.Range("L19").Value = -13.287904599659
$ws.Range("M19").Value = -11.439756415833
$ws.Range("N19").Value = -78.173241852487
# Row 20
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C20").Value = 2
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 1
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E20").Value = 100
$ws.Range("F20").NumberFormat = "#,##0"
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 66
$ws.Range("K20").Value = -28.787878787878
$ws.Range("L20").Value = -29.850746268656
$ws.Range("M20").Value = 104.347826086957
$ws.Range("N20").Value = -86.760563380281
# Row 21
$ws.Range("C21").Value = 62
$ws.Range("D21").Value = 76
$ws.Range("E21").Value = -18.421052631578
$ws.Range("F21").Value = 210
$ws.Range("G21").Value = 303
$ws.Range("H21").Value = -30.69306930693
$ws.Range("I21").Value = 3391
$ws.Range("J21").Value = 3785
$ws.Range("K21").Value = -10.409511228533
$ws.Range("L21").Value = -20.042442820089
$ws.Range("M21").Value = 11.399474375821
$ws.Range("N21").Value = -78.212541763042
# Row 22
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 100
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 165
$ws.Range("J22").Value = 198
$ws.Range("K22").Value = -16.666666666666
$ws.Range("L22").Value = -12.698412698412
$ws.Range("M22").Value = 11.486486486486
# Row 24
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 57
$ws.Range("E24").Value = -3.508771929824
$ws.Range("F24").Value = 302
$ws.Range("G24").Value = 263
$ws.Range("H24").Value = 14.828897338403
$ws.Range("I24").Value = 4314
$ws.Range("J24").Value = 3984
$ws.Range("K24").Value = 8.28313253012
$ws.Range("L24").Value = 27.144120247568
$ws.Range("M24").Value = -10.460772104607
# Row 25
$ws.Range("C25").Value = 53
$ws.Range("D25").Value = 49
$ws.Range("E25").Value = 8.163265306122
$ws.Range("F25").Value = 257
$ws.Range("G25").Value = 234
$ws.Range("H25").Value = 9.829059829059
$ws.Range("I25").Value = 3790
$ws.Range("J25").Value = 3572
$ws.Range("K25").Value = 6.103023516237
$ws.Range("L25").Value = 21.591273660571
# Row 26
$ws.Range("C26").Value = 22
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = 22.222222222222
$ws.Range("G26").Value = 89
$ws.Range("H26").Value = -4.494382022471
$ws.Range("I26").Value = 1052
$ws.Range("J26").Value = 1100
$ws.Range("K26").Value = -4.363636363636
$ws.Range("L26").Value = 14.09978308026
$ws.Range("M26").Value = 71.615008156606
# Row 27
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 36
$ws.Range("K27").Value = 38.461538461538
$ws.Range("L27").Value = 9.090909090909
# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = -6.666666666666
$ws.Range("I28").Value = 229
$ws.Range("J28").Value = 214
$ws.Range("K28").Value = 7.009345794392
$ws.Range("L28").Value = 3.619909502262
# Row 31
$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = -75
$ws.Range("I31").Value = 27
$ws.Range("J31").Value = 24
$ws.Range("K31").Value = 12.5
$ws.Range("L31").Value = 22.727272727272
